$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs)
$ws.Range("G2").Value = 29.358248
$ws.Range("H2").Value = 58.716496
$ws.Range("I2").Value = 0.02696851982721014
$ws.Range("J2").Value = 0.01839353552869724
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0108025
$ws.Range("N2").Value = 0.021605
$ws.Range("Q2").Value = 0.31714247402
$ws.Range("R2").Value = 1.26856989608
$ws.Range("S2").Value = 0.02696851982721014
$ws.Range("T2").Value = 0.01839353552869724

# Row 3 (FAPs)
$ws.Range("I3").Value = 0.9319510118584139
$ws.Range("J3").Value = 0.9534379801405131
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.0108025
$ws.Range("N3").Value = 0.021605
$ws.Range("Q3").Value = 10.9594909716925
$ws.Range("R3").Value = 65.756945830155
$ws.Range("S3").Value = 0.9319510118584139
$ws.Range("T3").Value = 0.9534379801405131

# Row 4 (Inflammatory-Mac)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.105487
$ws.Range("H4").Value = 0.316461
$ws.Range("I4").Value = 0.0000969004775425603
$ws.Range("J4").Value = 0.00009913460515332959
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.0108025
$ws.Range("N4").Value = 0.021605
$ws.Range("Q4").Value = 0.0011395233175
$ws.Range("R4").Value = 0.006837139904999999
$ws.Range("S4").Value = 0.0000969004775425603
$ws.Range("T4").Value = 0.00009913460515332959

# Row 5 (MuSCs)
$ws.Range("G5").Value = 44.2416095
$ws.Range("H5").Value = 88.48321900000001
$ws.Range("I5").Value = 0.040640392539379
$ws.Range("J5").Value = 0.02771826221322878
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.0108025
$ws.Range("N5").Value = 0.021605
$ws.Range("Q5").Value = 0.47791998662375
$ws.Range("R5").Value = 1.911679946495
$ws.Range("S5").Value = 0.040640392539379
$ws.Range("T5").Value = 0.02771826221322878

# Row 6 (Neutrophils)
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1512516666666667
$ws.Range("H6").Value = 0.453755
$ws.Range("I6").Value = 0.0001389399521183478
$ws.Range("J6").Value = 0.0001421433376035248
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.0108025
$ws.Range("N6").Value = 0.021605
$ws.Range("Q6").Value = 0.001633896129166667
$ws.Range("R6").Value = 0.009803376775000001
$ws.Range("S6").Value = 0.0001389399521183478
$ws.Range("T6").Value = 0.0001421433376035248

# Row 7 (Resolving-Mac)
$ws.Range("G7").Value = 0.222333
$ws.Range("H7").Value = 0.666999
$ws.Range("I7").Value = 0.0002042353453361084
$ws.Range("J7").Value = 0.0002089441748040539
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.0108025
$ws.Range("N7").Value = 0.021605
$ws.Range("Q7").Value = 0.0024017522325
$ws.Range("R7").Value = 0.014410513395
$ws.Range("S7").Value = 0.0002042353453361084
$ws.Range("T7").Value = 0.0002089441748040539
